# Generate Report for Handback
# Refresh the timestamp columns on the Overview / per-locale sheets to
# reflect the latest handoff/handback generation run (new xliff batch
# for bce2a8a9-9065-499c-bf50-4df3092b2d37.md).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for bce2a8a9-...md (row 2)
$wsOverview.Range("G2").Value = "2016-08-30 19:20:29"

# zh-cn sheet (row 2 = bce2a8a9-...md)
$wsZhCn.Range("H2").Value = "2016-08-30 19:20:23"   # Correspond Handoff Datetime
$wsZhCn.Range("K2").Value = "2016-08-30 19:20:41"   # Correspond Handback DateTime

# de-de sheet (row 2 = bce2a8a9-...md)
$wsDeDe.Range("H2").Value = "2016-08-30 19:20:29"   # Correspond Handoff Datetime (same generate run as Overview)
$wsDeDe.Range("K2").Value = "2016-08-30 19:20:49"   # Correspond Handback DateTime
